$d = $word.ActiveDocument

for ($i = 1; $i -le 5; $i++) {
    $old = "<id>p154r_a$i</id>"
    $new = "<id>p154r_$i</id>"
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}
